$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "fontsize" column header in S1
$ws.Range("S1").Value = "fontsize"

# Add fontsize values for each data row
$ws.Range("S2").Value = 20
$ws.Range("S3").Value = 16
$ws.Range("S4").Value = 18

# Update selection to match the post-edit state
$ws.Range("S1").Select()
